$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting so numeric-looking
# values (e.g. "583.38") are not coerced into actual numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.593.27"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").Value = "3.629.36"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D5").Value = "583.38"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("D6").Value = "175.65"
$ws.Range("E6").Value = "  -4.15%  "
$ws.Range("D7").Value = "3.621.77"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("D11").Value = "6.85"
$ws.Range("E11").Value = "  +16.59%  "
$ws.Range("D13").Value = "48.41"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "4.215.84"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "673.90"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "3.628.56"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "70.675.33"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  -4.50%  "
$ws.Range("D22").Value = "11.49"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "0.938"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "17.13"
$ws.Range("E24").Value = "  -4.33%  "
$ws.Range("D25").Value = "99.87"
$ws.Range("E25").Value = "  -5.15%  "
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "9.86"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "34.62"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").Value = "9.08"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  -5.65%  "
$ws.Range("D33").Value = "7.58"
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("E34").Value = "  -6.58%  "
$ws.Range("D35").Value = "3.97"
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("D36").Value = "578.08"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").Value = "11.09"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").Value = "58.46"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").Value = "3.568.12"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").Value = "0.346"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "34.46"
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("E46").Value = "  -6.70%  "
$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "137.03"
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").Value = "2.90"
$ws.Range("E51").Value = "  -2.86%  "
